# ============================================================================
# Edit: add "2022-Q3" quarter data to the 金力永磁 (300748) holdings workbook.
#
#   1. Insert a new "2022-Q3" row at the top of the per-quarter totals on
#      the "总计" sheet (pushing the existing quarters down by one row).
#   2. Insert a brand-new worksheet named "2022-Q3" right after "总计" so
#      the tab order becomes: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4,
#      2021-Q3, 2021-Q2, 2021-Q1 - and populate it with the per-fund
#      holding detail published for that quarter.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# 1. "总计" (summary) sheet: shift the existing 6 quarters down one row and
#    write the new 2022-Q3 totals into row 2.
# ----------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryRows = @(
    @("2022-Q3", 28, 2.05),
    @("2022-Q2", 21, 4.6),
    @("2022-Q1", 15, 2.74),
    @("2021-Q4", 10, 4.74),
    @("2021-Q3", 2, 1.25),
    @("2021-Q2", 8, 2.78),
    @("2021-Q1", 2, 0.05)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $summaryRows[$i][0]
    $summary.Cells.Item($r, 3).Value = $summaryRows[$i][1]
    $summary.Cells.Item($r, 4).Value = $summaryRows[$i][2]
}

# Rows 2-7 already existed (the sheet used to run from row 2 to row 7), so
# their "A" index cell already has the bold/centered style used throughout
# column A. Row 8 is brand new (the sheet used to stop at row 7) - clone
# the index-column style from the row directly above it.
$summary.Cells.Item(7, 1).Copy()
$summary.Cells.Item(8, 1).PasteSpecial(-4122)

# ----------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ----------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item(1))
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# Seed the header-row style (bold font + border, centered) and the
# index-column style by copying the already-styled cells from "总计" -
# this reuses the exact same style entries the other six quarter sheets
# use instead of minting new ones.
$summary.Cells.Item(1, 2).Copy()
$q3.Cells.Item(1, 2).PasteSpecial(-4122)
$summary.Cells.Item(1, 1).Copy()
$q3.Cells.Item(2, 1).PasteSpecial(-4122)

for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}
for ($c = 3; $c -le 8; $c++) {
    $q3.Cells.Item(1, 2).Copy()
    $q3.Cells.Item(1, $c).PasteSpecial(-4122)
}
for ($r = 3; $r -le 29; $r++) {
    $q3.Cells.Item(2, 1).Copy()
    $q3.Cells.Item($r, 1).PasteSpecial(-4122)
}

# Per-fund detail rows for 2022-Q3 - 基金代码/基金规模/股票总仓位/仓位占比/
# 持有市值 are published as plain text (leading "'" forces text so Excel
# doesn't coerce the numeric-looking strings into numbers), 仓位排名 is a
# real integer, matching how every other quarter sheet stores this table.
$data = @(
    @("'470009", "汇添富民营活力混合A", "'24.14", "'90.04", "'3.07", "'0.7411", 8),
    @("'690011", "民生加银积极成长混合", "'4.51", "'92.56", "'7.48", "'0.3373", 4),
    @("'005314", "万家中证1000指数增强C", "'14.28", "'94.11", "'1.12", "'0.1599", 2),
    @("'000884", "民生加银优选股票", "'1.99", "'92.47", "'7.50", "'0.1492", 4),
    @("'005313", "万家中证1000指数增强A", "'13.25", "'94.11", "'1.12", "'0.1484", 2),
    @("'001197", "长盛转型升级主题灵活配置混合", "'3.51", "'82.07", "'2.95", "'0.1035", 6),
    @("'005457", "景顺长城量化小盘股票", "'6.57", "'93.58", "'1.54", "'0.1012", 6),
    @("'013296", "民生加银聚优精选混合", "'1.26", "'92.34", "'7.53", "'0.0949", 3),
    @("'002210", "创金合信量化多因子股票A", "'2.39", "'91.71", "'1.30", "'0.0311", 2),
    @("'011888", "民生加银周期优选混合型证券投资基金A", "'0.39", "'92.24", "'7.59", "'0.0296", 3),
    @("'015496", "景顺中证1000指数增强C", "'1.83", "'92.63", "'1.60", "'0.0293", 7),
    @("'229002", "泰达宏利逆向策略混合", "'1.59", "'91.90", "'1.65", "'0.0262", 7),
    @("'001017", "泰达宏利改革动力量化策略灵活配置混合A", "'1.14", "'91.83", "'2.18", "'0.0249", 4),
    @("'009128", "明亚价值长青混合A", "'0.38", "'57.73", "'3.67", "'0.0139", 4),
    @("'015495", "景顺中证1000指数增强A", "'0.69", "'92.63", "'1.60", "'0.0110", 7),
    @("'003865", "创金合信量化多因子股票C", "'0.75", "'91.71", "'1.30", "'0.0098", 2),
    @("'004730", "建信量化事件驱动股票", "'0.46", "'81.67", "'1.90", "'0.0087", 5),
    @("'003647", "创金合信中证1000指数增强C", "'0.68", "'90.54", "'1.25", "'0.0085", 6),
    @("'011889", "民生加银周期优选混合型证券投资基金C", "'0.10", "'92.24", "'7.59", "'0.0076", 3),
    @("'320016", "诺安多策略混合", "'0.17", "'76.95", "'3.81", "'0.0065", 8),
    @("'003646", "创金合信中证1000指数增强A", "'0.35", "'90.54", "'1.25", "'0.0044", 6),
    @("'009514", "创金合信同顺创业板精选股票C", "'0.16", "'92.10", "'2.49", "'0.0040", 4),
    @("'009513", "创金合信同顺创业板精选股票A", "'0.09", "'92.10", "'2.49", "'0.0022", 4),
    @("'001419", "泰达宏利新思路灵活配置混合A", "'0.09", "'33.42", "'0.87", "'0.0008", 6),
    @("'003550", "泰达宏利改革动力量化策略灵活配置混合C", "'0.01", "'91.83", "'2.18", "'0.0002", 4),
    @("'960014", "汇添富民营活力混合 O", "'0.00", "'90.04", "'3.07", 0, 8),
    @("'002314", "泰达宏利新思路灵活配置混合B", "'0.00", "'33.42", "'0.87", 0, 6),
    @("'009129", "明亚价值长青混合C", "'0.00", "'57.73", "'3.67", 0, 4)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $q3.Cells.Item($r, 1).Value = $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $q3.Cells.Item($r, $c + 2).Value = $row[$c]
    }
}

# Keep "总计" as the active/selected sheet, matching the workbook's
# original view state (Worksheets.Add() activates the newly inserted
# sheet by default).
$summary.Activate()
